$d = $word.ActiveDocument

# 1) "=> Aumentamos X cuando se hace un lanzamiento..." -> "... se emite un lanzamiento..."
$d.Content.Find.Execute(
    "se hace un lanzamiento completo a",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "se emite un lanzamiento completo a", 2) | Out-Null

# 2) "=> Aumentamos Y cuando se lanza una corrección a la versión X." ->
#    "=> Aumentamos Y cuando se emite una corrección sobre la última versión X. en producción"
$d.Content.Find.Execute(
    "se lanza una corrección a la versión X.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "se emite una corrección sobre la última versión X. en producción", 2) | Out-Null

# 3) Drop the trailing parenthetical note about abbreviations
$d.Content.Find.Execute(
    "X.Y.1.0, X.Y.2.0, X.Y.3.0, etc. (nota: se pueden abreviar a X.Y.1, X.Y.2, X.Y.3, etc.)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "X.Y.1.0, X.Y.2.0, X.Y.3.0, etc.", 2) | Out-Null

# 4) "emite una versión para pruebas" -> "emite una versión "para pruebas" internas"
$d.Content.Find.Execute(
    "emite una versión para pruebas",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "emite una versión “para pruebas” internas", 2) | Out-Null

# 5) "empezamos a trabajar después de esa emisión." -> "... después de cualquier emisión."
$d.Content.Find.Execute(
    "empezamos a trabajar después de esa emisión.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "empezamos a trabajar después de cualquier emisión.", 2) | Out-Null
